# Swap the contents of columns D (codeforiati:group-name) and E (codeforiati:group-code)
# for every used row, including the header row, so that column D becomes the group
# code and column E becomes the group name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()

    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
